$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.697207333333333
$ws.Range("H2").Value = 29.091622
$ws.Range("I2").Value = 0.3243108558382838
$ws.Range("J2").Value = 0.3243108558382838
$ws.Range("M2").Value = 21.85308466666666
$ws.Range("N2").Value = 65.559254
$ws.Range("O2").Value = 0.407053040353553
$ws.Range("P2").Value = 0.407053040353553
$ws.Range("Q2").Value = 211.9138928855542
$ws.Range("R2").Value = 1907.225035969988
$ws.Range("S2").Value = 0.1320117198886362
$ws.Range("T2").Value = 0.1320117198886362
$ws.Range("G3").Value = 9.697207333333333
$ws.Range("H3").Value = 29.091622
$ws.Range("I3").Value = 0.3243108558382838
$ws.Range("J3").Value = 0.3243108558382838
$ws.Range("O3").Value = 0.1342711086924142
$ws.Range("P3").Value = 0.1342711086924142
$ws.Range("Q3").Value = 69.90222532266245
$ws.Range("R3").Value = 629.1200279039621
$ws.Range("S3").Value = 0.04354557817439209
$ws.Range("T3").Value = 0.04354557817439208
$ws.Range("G4").Value = 9.697207333333333
$ws.Range("H4").Value = 29.091622
$ws.Range("I4").Value = 0.3243108558382838
$ws.Range("J4").Value = 0.3243108558382838
$ws.Range("M4").Value = 11.375406
$ws.Range("N4").Value = 34.126218
$ws.Range("O4").Value = 0.2118874139822907
$ws.Range("P4").Value = 0.2118874139822907
$ws.Range("Q4").Value = 110.309670482844
$ws.Range("R4").Value = 992.7870343455961
$ws.Range("S4").Value = 0.06871738856995745
$ws.Range("T4").Value = 0.06871738856995743
$ws.Range("G5").Value = 9.697207333333333
$ws.Range("H5").Value = 29.091622
$ws.Range("I5").Value = 0.3243108558382838
$ws.Range("J5").Value = 0.3243108558382838
$ws.Range("M5").Value = 3.401340666666667
$ws.Range("N5").Value = 10.204022
$ws.Range("O5").Value = 0.06335609277882483
$ws.Range("P5").Value = 0.06335609277882483
$ws.Range("Q5").Value = 32.98350565596489
$ws.Range("R5").Value = 296.851550903684
$ws.Range("S5").Value = 0.02054706867167039
$ws.Range("T5").Value = 0.02054706867167039
$ws.Range("G6").Value = 9.697207333333333
$ws.Range("H6").Value = 29.091622
$ws.Range("I6").Value = 0.3243108558382838
$ws.Range("J6").Value = 0.3243108558382838
$ws.Range("M6").Value = 9.847764666666666
$ws.Range("N6").Value = 29.543294
$ws.Range("O6").Value = 0.1834323441929172
$ws.Range("P6").Value = 0.1834323441929172
$ws.Range("Q6").Value = 95.49581574254088
$ws.Range("R6").Value = 859.462341682868
$ws.Range("S6").Value = 0.05948910053362762
$ws.Range("T6").Value = 0.0594891005336276
$ws.Range("I7").Value = 0.2826325233457075
$ws.Range("J7").Value = 0.2826325233457074
$ws.Range("M7").Value = 21.85308466666666
$ws.Range("N7").Value = 65.559254
$ws.Range("O7").Value = 0.407053040353553
$ws.Range("P7").Value = 0.407053040353553
$ws.Range("Q7").Value = 184.68009072173
$ws.Range("R7").Value = 1662.12081649557
$ws.Range("S7").Value = 0.1150464279306668
$ws.Range("T7").Value = 0.1150464279306668
$ws.Range("I8").Value = 0.2826325233457075
$ws.Range("J8").Value = 0.2826325233457074
$ws.Range("O8").Value = 0.1342711086924142
$ws.Range("P8").Value = 0.1342711086924142
$ws.Range("S8").Value = 0.0379493822621628
$ws.Range("T8").Value = 0.03794938226216278
$ws.Range("I9").Value = 0.2826325233457075
$ws.Range("J9").Value = 0.2826325233457074
$ws.Range("M9").Value = 11.375406
$ws.Range("N9").Value = 34.126218
$ws.Range("O9").Value = 0.2118874139822907
$ws.Range("P9").Value = 0.2118874139822907
$ws.Range("Q9").Value = 96.13338547491001
$ws.Range("R9").Value = 865.2004692741901
$ws.Range("S9").Value = 0.05988627447901138
$ws.Range("T9").Value = 0.05988627447901135
$ws.Range("I10").Value = 0.2826325233457075
$ws.Range("J10").Value = 0.2826325233457074
$ws.Range("M10").Value = 3.401340666666667
$ws.Range("N10").Value = 10.204022
$ws.Range("O10").Value = 0.06335609277882483
$ws.Range("P10").Value = 0.06335609277882483
$ws.Range("Q10").Value = 28.74467895389
$ws.Range("R10").Value = 258.70211058501
$ws.Range("S10").Value = 0.01790649237140402
$ws.Range("T10").Value = 0.01790649237140401
$ws.Range("I11").Value = 0.2826325233457075
$ws.Range("J11").Value = 0.2826325233457074
$ws.Range("M11").Value = 9.847764666666666
$ws.Range("N11").Value = 29.543294
$ws.Range("O11").Value = 0.1834323441929172
$ws.Range("P11").Value = 0.1834323441929172
$ws.Range("Q11").Value = 83.22331148153
$ws.Range("R11").Value = 749.00980333377
$ws.Range("S11").Value = 0.05184394630246252
$ws.Range("T11").Value = 0.05184394630246249
$ws.Range("G12").Value = 3.910524
$ws.Range("H12").Value = 11.731572
$ws.Range("I12").Value = 0.1307825378608469
$ws.Range("J12").Value = 0.1307825378608469
$ws.Range("M12").Value = 21.85308466666666
$ws.Range("N12").Value = 65.559254
$ws.Range("O12").Value = 0.407053040353553
$ws.Range("P12").Value = 0.407053040353553
$ws.Range("Q12").Value = 85.457012063032
$ws.Range("R12").Value = 769.1131085672879
$ws.Range("S12").Value = 0.0532354296614114
$ws.Range("T12").Value = 0.05323542966141138
$ws.Range("G13").Value = 3.910524
$ws.Range("H13").Value = 11.731572
$ws.Range("I13").Value = 0.1307825378608469
$ws.Range("J13").Value = 0.1307825378608469
$ws.Range("O13").Value = 0.1342711086924142
$ws.Range("P13").Value = 0.1342711086924142
$ws.Range("Q13").Value = 28.188974452268
$ws.Range("R13").Value = 253.700770070412
$ws.Range("S13").Value = 0.01756031635618355
$ws.Range("T13").Value = 0.01756031635618355
$ws.Range("G14").Value = 3.910524
$ws.Range("H14").Value = 11.731572
$ws.Range("I14").Value = 0.1307825378608469
$ws.Range("J14").Value = 0.1307825378608469
$ws.Range("M14").Value = 11.375406
$ws.Range("N14").Value = 34.126218
$ws.Range("O14").Value = 0.2118874139822907
$ws.Range("P14").Value = 0.2118874139822907
$ws.Range("Q14").Value = 44.483798172744
$ws.Range("R14").Value = 400.354183554696
$ws.Range("S14").Value = 0.02771117374137588
$ws.Range("T14").Value = 0.02771117374137587
$ws.Range("G15").Value = 3.910524
$ws.Range("H15").Value = 11.731572
$ws.Range("I15").Value = 0.1307825378608469
$ws.Range("J15").Value = 0.1307825378608469
$ws.Range("M15").Value = 3.401340666666667
$ws.Range("N15").Value = 10.204022
$ws.Range("O15").Value = 0.06335609277882483
$ws.Range("P15").Value = 0.06335609277882483
$ws.Range("Q15").Value = 13.301024309176
$ws.Range("R15").Value = 119.709218782584
$ws.Range("S15").Value = 0.008285870602561989
$ws.Range("T15").Value = 0.008285870602561987
$ws.Range("G16").Value = 3.910524
$ws.Range("H16").Value = 11.731572
$ws.Range("I16").Value = 0.1307825378608469
$ws.Range("J16").Value = 0.1307825378608469
$ws.Range("M16").Value = 9.847764666666666
$ws.Range("N16").Value = 29.543294
$ws.Range("O16").Value = 0.1834323441929172
$ws.Range("P16").Value = 0.1834323441929172
$ws.Range("Q16").Value = 38.509920075352
$ws.Range("R16").Value = 346.589280678168
$ws.Range("S16").Value = 0.02398974749931409
$ws.Range("T16").Value = 0.02398974749931409
$ws.Range("G17").Value = 3.066674
$ws.Range("H17").Value = 9.200022000000001
$ws.Range("I17").Value = 0.1025610400324547
$ws.Range("J17").Value = 0.1025610400324547
$ws.Range("M17").Value = 21.85308466666666
$ws.Range("N17").Value = 65.559254
$ws.Range("O17").Value = 0.407053040353553
$ws.Range("P17").Value = 0.407053040353553
$ws.Range("Q17").Value = 67.01628656706534
$ws.Range("R17").Value = 603.146579103588
$ws.Range("S17").Value = 0.04174778316703315
$ws.Range("T17").Value = 0.04174778316703314
$ws.Range("G18").Value = 3.066674
$ws.Range("H18").Value = 9.200022000000001
$ws.Range("I18").Value = 0.1025610400324547
$ws.Range("J18").Value = 0.1025610400324547
$ws.Range("O18").Value = 0.1342711086924142
$ws.Range("P18").Value = 0.1342711086924142
$ws.Range("Q18").Value = 22.10608988448467
$ws.Range("R18").Value = 198.954808960362
$ws.Range("S18").Value = 0.01377098455380477
$ws.Range("T18").Value = 0.01377098455380477
$ws.Range("G19").Value = 3.066674
$ws.Range("H19").Value = 9.200022000000001
$ws.Range("I19").Value = 0.1025610400324547
$ws.Range("J19").Value = 0.1025610400324547
$ws.Range("M19").Value = 11.375406
$ws.Range("N19").Value = 34.126218
$ws.Range("O19").Value = 0.2118874139822907
$ws.Range("P19").Value = 0.2118874139822907
$ws.Range("Q19").Value = 34.884661819644
$ws.Range("R19").Value = 313.961956376796
$ws.Range("S19").Value = 0.02173139354781102
$ws.Range("T19").Value = 0.02173139354781101
$ws.Range("G20").Value = 3.066674
$ws.Range("H20").Value = 9.200022000000001
$ws.Range("I20").Value = 0.1025610400324547
$ws.Range("J20").Value = 0.1025610400324547
$ws.Range("M20").Value = 3.401340666666667
$ws.Range("N20").Value = 10.204022
$ws.Range("O20").Value = 0.06335609277882483
$ws.Range("P20").Value = 0.06335609277882483
$ws.Range("Q20").Value = 10.43080298760933
$ws.Range("R20").Value = 93.877226888484
$ws.Range("S20").Value = 0.006497866767788968
$ws.Range("T20").Value = 0.006497866767788966
$ws.Range("G21").Value = 3.066674
$ws.Range("H21").Value = 9.200022000000001
$ws.Range("I21").Value = 0.1025610400324547
$ws.Range("J21").Value = 0.1025610400324547
$ws.Range("M21").Value = 9.847764666666666
$ws.Range("N21").Value = 29.543294
$ws.Range("O21").Value = 0.1834323441929172
$ws.Range("P21").Value = 0.1834323441929172
$ws.Range("Q21").Value = 30.19988386138534
$ws.Range("R21").Value = 271.798954752468
$ws.Range("S21").Value = 0.01881301199601679
$ws.Range("T21").Value = 0.01881301199601678
$ws.Range("G22").Value = 4.775574
$ws.Range("H22").Value = 14.326722
$ws.Range("I22").Value = 0.1597130429227071
$ws.Range("J22").Value = 0.159713042922707
$ws.Range("M22").Value = 21.85308466666666
$ws.Range("N22").Value = 65.559254
$ws.Range("O22").Value = 0.407053040353553
$ws.Range("P22").Value = 0.407053040353553
$ws.Range("Q22").Value = 104.361022953932
$ws.Range("R22").Value = 939.2492065853879
$ws.Range("S22").Value = 0.06501167970580543
$ws.Range("T22").Value = 0.06501167970580542
$ws.Range("G23").Value = 4.775574
$ws.Range("H23").Value = 14.326722
$ws.Range("I23").Value = 0.1597130429227071
$ws.Range("J23").Value = 0.159713042922707
$ws.Range("O23").Value = 0.1342711086924142
$ws.Range("P23").Value = 0.1342711086924142
$ws.Range("Q23").Value = 34.424679015118
$ws.Range("R23").Value = 309.822111136062
$ws.Range("S23").Value = 0.02144484734587102
$ws.Range("T23").Value = 0.02144484734587102
$ws.Range("G24").Value = 4.775574
$ws.Range("H24").Value = 14.326722
$ws.Range("I24").Value = 0.1597130429227071
$ws.Range("J24").Value = 0.159713042922707
$ws.Range("M24").Value = 11.375406
$ws.Range("N24").Value = 34.126218
$ws.Range("O24").Value = 0.2118874139822907
$ws.Range("P24").Value = 0.2118874139822907
$ws.Range("Q24").Value = 54.32409313304399
$ws.Range("R24").Value = 488.916838197396
$ws.Range("S24").Value = 0.033841183644135
$ws.Range("T24").Value = 0.03384118364413499
$ws.Range("G25").Value = 4.775574
$ws.Range("H25").Value = 14.326722
$ws.Range("I25").Value = 0.1597130429227071
$ws.Range("J25").Value = 0.159713042922707
$ws.Range("M25").Value = 3.401340666666667
$ws.Range("N25").Value = 10.204022
$ws.Range("O25").Value = 0.06335609277882483
$ws.Range("P25").Value = 0.06335609277882483
$ws.Range("Q25").Value = 16.243354052876
$ws.Range("R25").Value = 146.190186475884
$ws.Range("S25").Value = 0.01011879436539946
$ws.Range("T25").Value = 0.01011879436539946
$ws.Range("G26").Value = 4.775574
$ws.Range("H26").Value = 14.326722
$ws.Range("I26").Value = 0.1597130429227071
$ws.Range("J26").Value = 0.159713042922707
$ws.Range("M26").Value = 9.847764666666666
$ws.Range("N26").Value = 29.543294
$ws.Range("O26").Value = 0.1834323441929172
$ws.Range("P26").Value = 0.1834323441929172
$ws.Range("Q26").Value = 47.028728900252
$ws.Range("R26").Value = 423.258560102268
$ws.Range("S26").Value = 0.02929653786149616
$ws.Range("T26").Value = 0.02929653786149615
